$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Cells.Item(2, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "30.145.86"
$cD.Style = $origStyle
$ws.Cells.Item(2, 5).Value = "  -3.37%  "

$cD = $ws.Cells.Item(3, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.862.87"
$cD.Style = $origStyle
$ws.Cells.Item(3, 5).Value = "  -4.00%  "

$cD = $ws.Cells.Item(4, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.000"
$cD.Style = $origStyle
$ws.Cells.Item(4, 5).Value = "  +0.01%  "

$cD = $ws.Cells.Item(5, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "233.65"
$cD.Style = $origStyle
$ws.Cells.Item(5, 5).Value = "  -3.48%  "

$cD = $ws.Cells.Item(6, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.000"
$cD.Style = $origStyle
$ws.Cells.Item(6, 5).Value = "  +0.00%  "

$cD = $ws.Cells.Item(7, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.4661"
$cD.Style = $origStyle
$ws.Cells.Item(7, 5).Value = "  -2.93%  "

$cD = $ws.Cells.Item(8, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.2831"
$cD.Style = $origStyle
$ws.Cells.Item(8, 5).Value = "  -2.79%  "

$cD = $ws.Cells.Item(9, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.06545"
$cD.Style = $origStyle
$ws.Cells.Item(9, 5).Value = "  -3.52%  "

$cD = $ws.Cells.Item(10, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "20.13"
$cD.Style = $origStyle
$ws.Cells.Item(10, 5).Value = "  -0.49%  "

$cD = $ws.Cells.Item(11, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.07817"
$cD.Style = $origStyle
$ws.Cells.Item(11, 5).Value = "  -0.32%  "

$cD = $ws.Cells.Item(12, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "96.09"
$cD.Style = $origStyle
$ws.Cells.Item(12, 5).Value = "  -7.86%  "

$cD = $ws.Cells.Item(13, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.856.09"
$cD.Style = $origStyle
$ws.Cells.Item(13, 5).Value = "  -4.46%  "

$cD = $ws.Cells.Item(14, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "5.118"
$cD.Style = $origStyle
$ws.Cells.Item(14, 5).Value = "  -3.51%  "

$cD = $ws.Cells.Item(15, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.6687"
$cD.Style = $origStyle
$ws.Cells.Item(15, 5).Value = "  -4.34%  "

$cD = $ws.Cells.Item(16, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "280.06"
$cD.Style = $origStyle
$ws.Cells.Item(16, 5).Value = "  -5.19%  "

$cD = $ws.Cells.Item(17, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "30.168.34"
$cD.Style = $origStyle
$ws.Cells.Item(17, 5).Value = "  -3.30%  "

$cD = $ws.Cells.Item(18, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.000"
$cD.Style = $origStyle
$ws.Cells.Item(18, 5).Value = "  +0.04%  "

$cD = $ws.Cells.Item(19, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "5.449"
$cD.Style = $origStyle
$ws.Cells.Item(19, 5).Value = "  -2.22%  "

$cD = $ws.Cells.Item(20, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "12.61"
$cD.Style = $origStyle
$ws.Cells.Item(20, 5).Value = "  -3.08%  "

$cD = $ws.Cells.Item(21, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.098.38"
$cD.Style = $origStyle
$ws.Cells.Item(21, 5).Value = "  -5.01%  "

$cD = $ws.Cells.Item(22, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.000007228"
$cD.Style = $origStyle
$ws.Cells.Item(22, 5).Value = "  -5.03%  "

$ws.Cells.Item(23, 5).Value = "  +0.04%  "

$cD = $ws.Cells.Item(24, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "6.139"
$cD.Style = $origStyle
$ws.Cells.Item(24, 5).Value = "  -4.54%  "

$cD = $ws.Cells.Item(25, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "9.307"
$cD.Style = $origStyle
$ws.Cells.Item(25, 5).Value = "  -2.68%  "

$cD = $ws.Cells.Item(26, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "165.37"
$cD.Style = $origStyle
$ws.Cells.Item(26, 5).Value = "  -2.39%  "

$cD = $ws.Cells.Item(27, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "18.88"
$cD.Style = $origStyle
$ws.Cells.Item(27, 5).Value = "  -4.79%  "

$cD = $ws.Cells.Item(28, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.903"
$cD.Style = $origStyle
$ws.Cells.Item(28, 5).Value = "  -9.55%  "

$ws.Cells.Item(29, 5).Value = "  -3.51%  "

$cD = $ws.Cells.Item(30, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.09615"
$cD.Style = $origStyle
$ws.Cells.Item(30, 5).Value = "  -4.52%  "

$cD = $ws.Cells.Item(31, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "4.399"
$cD.Style = $origStyle
$ws.Cells.Item(31, 5).Value = "  -5.11%  "

$cD = $ws.Cells.Item(32, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.469"
$cD.Style = $origStyle
$ws.Cells.Item(32, 5).Value = "  -4.28%  "

$cD = $ws.Cells.Item(33, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "4.111"
$cD.Style = $origStyle
$ws.Cells.Item(33, 5).Value = "  -5.27%  "

$cD = $ws.Cells.Item(34, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.04658"
$cD.Style = $origStyle
$ws.Cells.Item(34, 5).Value = "  -3.75%  "

$cD = $ws.Cells.Item(35, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.7007"
$cD.Style = $origStyle
$ws.Cells.Item(35, 5).Value = "  -5.20%  "

$cD = $ws.Cells.Item(36, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.095"
$cD.Style = $origStyle
$ws.Cells.Item(36, 5).Value = "  -3.39%  "

$cD = $ws.Cells.Item(37, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.711"
$cD.Style = $origStyle
$ws.Cells.Item(37, 5).Value = "  -0.76%  "

$cD = $ws.Cells.Item(38, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.01851"
$cD.Style = $origStyle
$ws.Cells.Item(38, 5).Value = "  -5.46%  "

$cD = $ws.Cells.Item(39, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "6.275"
$cD.Style = $origStyle
$ws.Cells.Item(39, 5).Value = "  -8.97%  "

$cD = $ws.Cells.Item(40, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.520"
$cD.Style = $origStyle
$ws.Cells.Item(40, 5).Value = "  -4.29%  "

$cD = $ws.Cells.Item(41, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "72.43"
$cD.Style = $origStyle
$ws.Cells.Item(41, 5).Value = "  -5.54%  "

$cD = $ws.Cells.Item(42, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.8527"
$cD.Style = $origStyle
$ws.Cells.Item(42, 5).Value = "  -2.32%  "

$cD = $ws.Cells.Item(43, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.922"
$cD.Style = $origStyle
$ws.Cells.Item(43, 5).Value = "  -5.59%  "

$cD = $ws.Cells.Item(44, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.000"
$cD.Style = $origStyle
$ws.Cells.Item(44, 5).Value = "  +0.01%  "

$cD = $ws.Cells.Item(45, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.4154"
$cD.Style = $origStyle

$cD = $ws.Cells.Item(46, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "103.26"
$cD.Style = $origStyle
$ws.Cells.Item(46, 5).Value = "  -2.62%  "

$cD = $ws.Cells.Item(47, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "986.60"
$cD.Style = $origStyle
$ws.Cells.Item(47, 5).Value = "  -3.12%  "

$cD = $ws.Cells.Item(48, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "7.157"
$cD.Style = $origStyle
$ws.Cells.Item(48, 5).Value = "  -5.76%  "

$cD = $ws.Cells.Item(49, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "9.151"
$cD.Style = $origStyle
$ws.Cells.Item(49, 5).Value = "  -1.07%  "

$cD = $ws.Cells.Item(50, 4)
$origStyle = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "34.15"
$cD.Style = $origStyle
$ws.Cells.Item(50, 5).Value = "  -3.13%  "

$ws.Cells.Item(51, 5).Value = "  -5.95%  "

